$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph robustly (instead of
# hard-coding an index) by searching the document content for its text.
$findRange = $d.Content
$found = $findRange.Find.Execute("Docente(s) Responsável(eis)")
if (-not $found) {
    throw "Could not find 'Docente(s) Responsável(eis)' heading paragraph"
}
$headingPara = $findRange.Paragraphs(1)

# Insert a brand-new paragraph right after the heading; it will hold the
# bulleted list of professors.
$headingPara.Range.InsertParagraphAfter()
$listPara = $headingPara.Next()
$listPara.Style = "ListBullet"

# Each professor on its own line, separated by manual line breaks (w:br),
# matching the existing "Créditos-aula/Créditos-trabalho/..." list further up
# in the document.
$names = @(
    '144651 - Antonio Fernando Sartori',
    '3577649 - Carlos Angelo Nunes',
    '471420 - Carlos Antonio Reis Pereira Baptista',
    '519033 - Carlos Yujiro Shigue',
    '3586455 - Cassius Olivio Figueiredo Terra Ruchert',
    '5840897 - Clodoaldo Saron',
    '5840963 - Daniela Camargo Vernilli',
    '6495737 - Durval Rodrigues Junior',
    '1033242 - Fábio Herbst Florenzano',
    '5983729 - Fernando Vernilli Junior',
    '5009972 - Gilberto Carvalho Coelho',
    '984972 - Hugo Ricardo Zschommler Sandim',
    '1176388 - Luiz Tadeu Fernandes Eleno',
    '7459752 - Maria Ismenia Sodero Toledo Faria',
    '5840622 - Miguel Justino Ribeiro Barboza',
    '2166002 - Sandra Giacomin Schneider',
    '1922320 - Sebastiao Ribeiro',
    '5840793 - Sérgio Schneider'
)

$lineBreak = [string][char]11
$pos = $listPara.Range.Start

for ($i = 0; $i -lt $names.Length; $i++) {
    $segment = $names[$i]
    if ($i -lt $names.Length - 1) {
        $segment = $segment + $lineBreak
    }
    $insertRange = $d.Range($pos, $pos)
    $insertRange.InsertAfter($segment)
    $pos = $pos + $segment.Length
}
